$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Sostituzione_0" record currently sits on row 6, right below the four
# "Creazione_*" records (rows 2-5) and above "Aggiornamento_0" (row 7).
# This edit moves the "Sostituzione_0" record up to row 2 (right after the
# header row), pushing the "Creazione_0".."Creazione_3" records down by one
# row each. The "Aggiornamento_0" row (row 7) keeps its own values.

$cols = 5

# Snapshot every row's values (by column) before we start overwriting cells,
# so later writes never clobber data we still need to read.
$snapshot = @{}
for ($r = 2; $r -le 6; $r++) {
    $vals = @()
    for ($c = 1; $c -le $cols; $c++) {
        $vals += , ($ws.Cells.Item($r, $c).Text)
    }
    $snapshot[$r] = $vals
}

# Row 2 becomes the old row 6 ("Sostituzione_0").
for ($c = 1; $c -le $cols; $c++) {
    $ws.Cells.Item(2, $c).Value = $snapshot[6][$c - 1]
}

# Rows 3-6 become the old rows 2-5 ("Creazione_0".."Creazione_3"), each
# shifted down by one row.
for ($r = 3; $r -le 6; $r++) {
    $srcRow = $r - 1
    for ($c = 1; $c -le $cols; $c++) {
        $ws.Cells.Item($r, $c).Value = $snapshot[$srcRow][$c - 1]
    }
}
